$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21.
foreach ($r in 2..21) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
